# descw-1256 model returns report JSON with totals
#
# The report-template workbook's placeholder tokens are rewritten so the
# model iterates `d.report` (with a nested `.totals` collection for grand
# totals) instead of the old `d.repotsByProjectWithTotals` /
# `d.grandTotals` shapes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-group / project-group / grand-total iterator placeholders.
$ws.Range("B12").Value = "{#r=d.report[i]}"
$ws.Range("B13").Value = "{#r1=d.report[i+1]}"
$ws.Range("B15").Value = "{#p=d.report[i].projects[i]}"
$ws.Range("B16").Value = "{#p1=d.report[i].projects[i+1]}"
$ws.Range("B18").Value = "{#gt=d.totals[0]}"

# Leave the active selection on the cell the author ended up editing.
$null = $ws.Range("B12").Select()
